# Auto-generated update of crypto price/volume data, and Filecoin/Mantle row swap.
# Every data cell in columns B:E of this sheet is stored as text (inline string),
# even when the text looks like a plain number (e.g. "1.00", "528.35"). Writing a
# bare numeric-looking string through the Excel object model would normally be
# auto-converted to a real number, so we prefix every new value with a leading
# apostrophe. That is exactly how Excel's UI forces "store as text" behavior, and
# it reproduces the original inline-string cell content without altering the
# visible text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '58.235.29' },
    @{ Cell = 'E2'; Value = '  -0.83%  ' },
    @{ Cell = 'D3'; Value = '3.122.68' },
    @{ Cell = 'E3'; Value = '  +1.15%  ' },
    @{ Cell = 'E4'; Value = '  -0.02%  ' },
    @{ Cell = 'D5'; Value = '528.35' },
    @{ Cell = 'E5'; Value = '  +1.00%  ' },
    @{ Cell = 'D6'; Value = '142.59' },
    @{ Cell = 'E6'; Value = '  -0.68%  ' },
    @{ Cell = 'E7'; Value = '  +0.05%  ' },
    @{ Cell = 'D8'; Value = '3.120.25' },
    @{ Cell = 'D9'; Value = '0.446' },
    @{ Cell = 'E9'; Value = '  +1.70%  ' },
    @{ Cell = 'E10'; Value = '  -2.61%  ' },
    @{ Cell = 'E11'; Value = '  -0.20%  ' },
    @{ Cell = 'D12'; Value = '0.394' },
    @{ Cell = 'E12'; Value = '  +2.37%  ' },
    @{ Cell = 'D13'; Value = '3.657.64' },
    @{ Cell = 'E13'; Value = '  +1.08%  ' },
    @{ Cell = 'D14'; Value = '0.135' },
    @{ Cell = 'E14'; Value = '  +3.59%  ' },
    @{ Cell = 'D15'; Value = '25.66' },
    @{ Cell = 'E15'; Value = '  -4.25%  ' },
    @{ Cell = 'E16'; Value = '  -0.73%  ' },
    @{ Cell = 'D17'; Value = '58.262.10' },
    @{ Cell = 'E17'; Value = '  -0.81%  ' },
    @{ Cell = 'D18'; Value = '3.115.81' },
    @{ Cell = 'E18'; Value = '  +0.96%  ' },
    @{ Cell = 'E19'; Value = '  -0.57%  ' },
    @{ Cell = 'D20'; Value = '12.82' },
    @{ Cell = 'E20'; Value = '  -0.66%  ' },
    @{ Cell = 'D21'; Value = '7.99' },
    @{ Cell = 'E21'; Value = '  -1.80%  ' },
    @{ Cell = 'D22'; Value = '343.26' },
    @{ Cell = 'E22'; Value = '  +0.35%  ' },
    @{ Cell = 'D23'; Value = '1.00' },
    @{ Cell = 'E23'; Value = '  -0.02%  ' },
    @{ Cell = 'E24'; Value = '  +1.53%  ' },
    @{ Cell = 'D25'; Value = '67.65' },
    @{ Cell = 'E25'; Value = '  +2.88%  ' },
    @{ Cell = 'D26'; Value = '0.169' },
    @{ Cell = 'E26'; Value = '  -0.99%  ' },
    @{ Cell = 'D27'; Value = '1.00' },
    @{ Cell = 'E27'; Value = '  +0.02%  ' },
    @{ Cell = 'D28'; Value = '0.0₃0932' },
    @{ Cell = 'E28'; Value = '  +1.15%  ' },
    @{ Cell = 'E29'; Value = '  +0.04%  ' },
    @{ Cell = 'D30'; Value = '7.37' },
    @{ Cell = 'E30'; Value = '  +1.59%  ' },
    @{ Cell = 'E31'; Value = '  -3.52%  ' },
    @{ Cell = 'D32'; Value = '1.87' },
    @{ Cell = 'E32'; Value = '  +1.52%  ' },
    @{ Cell = 'D33'; Value = '21.07' },
    @{ Cell = 'E33'; Value = '  +0.17%  ' },
    @{ Cell = 'E34'; Value = '  -1.56%  ' },
    @{ Cell = 'D35'; Value = '158.57' },
    @{ Cell = 'E35'; Value = '  +2.62%  ' },
    @{ Cell = 'D36'; Value = '4.73' },
    @{ Cell = 'E36'; Value = '  +2.51%  ' },
    @{ Cell = 'D37'; Value = '6.23' },
    @{ Cell = 'E37'; Value = '  +1.69%  ' },
    @{ Cell = 'D38'; Value = '26.30' },
    @{ Cell = 'E38'; Value = '  -2.25%  ' },
    @{ Cell = 'E39'; Value = '  -3.93%  ' },
    @{ Cell = 'E40'; Value = '  +11.79%  ' },
    @{ Cell = 'D41'; Value = '0.0666' },
    @{ Cell = 'E41'; Value = '  -2.42%  ' },
    @{ Cell = 'B42'; Value = 'Mantle' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' },
    @{ Cell = 'D42'; Value = '0.695' },
    @{ Cell = 'E42'; Value = '  +4.58%  ' },
    @{ Cell = 'B43'; Value = 'Filecoin' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D43'; Value = '3.98' },
    @{ Cell = 'E43'; Value = '  +1.98%  ' },
    @{ Cell = 'D44'; Value = '3.161.39' },
    @{ Cell = 'E44'; Value = '  +1.00%  ' },
    @{ Cell = 'E45'; Value = '  -0.46%  ' },
    @{ Cell = 'D46'; Value = '1.00' },
    @{ Cell = 'E46'; Value = '  -0.05%  ' },
    @{ Cell = 'D47'; Value = '0.0263' },
    @{ Cell = 'D48'; Value = '2.278.06' },
    @{ Cell = 'E48'; Value = '  +0.22%  ' },
    @{ Cell = 'E49'; Value = '  +3.95%  ' },
    @{ Cell = 'E50'; Value = '  +1.63%  ' },
    @{ Cell = 'D51'; Value = '20.64' },
    @{ Cell = 'E51'; Value = '  -0.61%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
